# VIRGINIA_2016.xlsx cleanup:
#  1. Rename header row (A1:D1) to short machine-friendly column codes.
#  2. Title-case the small Spanish connector words ("de", "del", "la",
#     "las", "los", "el", "y") wherever they appear as standalone tokens
#     inside the state/municipality name columns (A, B).
#  3. Fix a floating point last-bit representation drift in the
#     "Porcentaje de Matrículas" column (D) for rows holding the
#     5/5156 ratio.
#  4. Drop the trailing footnote/metadata rows (rows 1032-1036) that sat
#     below the data table, shrinking the used range to A1:D1030.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row -----------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Title-case connector words in columns A and B -------------------
$stopwords = @("de", "del", "la", "las", "los", "el", "y")

function Contains-CaseSensitive($arr, $item) {
    foreach ($x in $arr) {
        if ($x.Equals($item)) {
            return $true
        }
    }
    return $false
}

function Fix-SpanishConnectors($s) {
    $parts = $s.Split(" ")
    $out = @()
    foreach ($p in $parts) {
        if (Contains-CaseSensitive $stopwords $p) {
            $out += $p.Substring(0, 1).ToUpper() + $p.Substring(1)
        } else {
            $out += $p
        }
    }
    return ($out -join " ")
}

for ($r = 2; $r -le 1030; $r++) {
    foreach ($c in @(1, 2)) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ne $null -and $v -is [string]) {
            $fixed = Fix-SpanishConnectors($v)
            $cell.Value = $fixed
        }
    }
}

# --- 3. Floating point drift fix for the 5/5156 ratio -------------------
$target = 0.0009697439875872769
$newval = 0.0009697439875872768
for ($r = 2; $r -le 1030; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $v = $cell.Value2
    if ($v -ne $null -and ($v -is [double] -or $v -is [int])) {
        if ([math]::Abs($v - $target) -lt 0.0000000000000000001) {
            $cell.Value = $newval
        }
    }
}

# --- 4. Drop trailing footnote rows -------------------------------------
$ws.Range("A1032:A1036").EntireRow.Delete()
